# New weekly price record for "Macroferia Regional de Talca - Piña".
# The sheet lists the most recent week first (right after the header),
# so this record is inserted as the new row 141, pushing every existing
# data row (old 141..168) down by one (new 142..169).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 141, shifting rows 141:168 down to 142:169.
$ws.Rows.Item(141).Insert()

# Populate the new row with this week's record.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 44505
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100108
$ws.Range("H141").Value = "Tropicales y subtropicales"
$ws.Range("I141").Value = 100108005
$ws.Range("J141").Value = "Piña"
$ws.Range("K141").Value = "Caramelo"
$ws.Range("L141").Value = "Segunda"
$ws.Range("M141").Value = 240
$ws.Range("N141").Value = 17000
$ws.Range("O141").Value = 17000
$ws.Range("P141").Value = 17000
$ws.Range("Q141").Value = '$/caja 14 unidades'
$ws.Range("R141").Value = "Ecuador"
$ws.Range("S141").Value = 1214
$ws.Range("T141").Value = 14
